$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 1274.2
$ws.Range("I5").Value = 443
$ws.Range("J5").Value = 2521
$ws.Range("K5").Value = 443
$ws.Range("L5").Value = 2521
$ws.Range("M5").Value = -328
$ws.Range("N5").Value = -2751
$ws.Range("H11").Value = 185.38461
$ws.Range("I11").Value = 185.38461
$ws.Range("K11").Value = 185.38461
$ws.Range("M11").Value = -45.38461000000001
$ws.Range("H12").Value = 474.3
$ws.Range("I12").Value = 549
$ws.Range("K12").Value = 549
$ws.Range("M12").Value = -379
$ws.Range("H33").Value = 479.3846
$ws.Range("I33").Value = 303.1
$ws.Range("K33").Value = 303.1
$ws.Range("M33").Value = -74.10000000000002
$ws.Range("H43").Value = 1251018.8
$ws.Range("I43").Value = 1133.3334
$ws.Range("K43").Value = 1133.3334
$ws.Range("M43").Value = -1064.3334
$ws.Range("H51").Value = 4707.9585
$ws.Range("J51").Value = 4756.174
$ws.Range("L51").Value = 4756.174
$ws.Range("N51").Value = -5724.174
$ws.Range("H100").Value = 3022.5
$ws.Range("J100").Value = 3998.4443
$ws.Range("L100").Value = 3998.4443
$ws.Range("N100").Value = -5080.4443
$ws.Range("H101").Value = 1524.75
$ws.Range("I101").Value = 299
$ws.Range("K101").Value = 897
$ws.Range("M101").Value = 725
$ws.Range("H132").Value = 5870.636
$ws.Range("I132").Value = 6357.7
$ws.Range("K132").Value = 19073.1
$ws.Range("M132").Value = -16543.1
$ws.Range("H138").Value = 5183.364
$ws.Range("J138").Value = 8074.077
$ws.Range("L138").Value = 24222.231
$ws.Range("N138").Value = -34502.231
$ws.Range("H141").Value = 15082.909
$ws.Range("I141").Value = 21150.8
$ws.Range("J141").Value = 10026.333
$ws.Range("K141").Value = 63452.39999999999
$ws.Range("L141").Value = 30078.999
$ws.Range("M141").Value = -58272.39999999999
$ws.Range("N141").Value = -40438.999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3798.6
$ws.Range("I2").Value = 3359.8
$ws.Range("J2").Value = 4237.4
$ws.Range("K2").Value = 3359.8
$ws.Range("L2").Value = 4237.4
$ws.Range("M2").Value = -3246.8
$ws.Range("N2").Value = -4463.4
$ws.Range("H45").Value = 2512
$ws.Range("I45").Value = 1826
$ws.Range("K45").Value = 1826
$ws.Range("M45").Value = -1449
$ws.Range("H46").Value = 7623.25
$ws.Range("I46").Value = 8997
$ws.Range("J46").Value = 7165.3335
$ws.Range("K46").Value = 8997
$ws.Range("L46").Value = 7165.3335
$ws.Range("M46").Value = -8678
$ws.Range("N46").Value = -7803.3335
$ws.Range("H116").Value = 3798.6
$ws.Range("I116").Value = 3359.8
$ws.Range("J116").Value = 4237.4
$ws.Range("K116").Value = 3359.8
$ws.Range("L116").Value = 4237.4
$ws.Range("M116").Value = -1065.8
$ws.Range("N116").Value = -8825.4
$ws.Range("H125").Value = 149499.33
$ws.Range("J125").Value = 149499.33
$ws.Range("L125").Value = 149499.33
$ws.Range("N125").Value = -159339.33

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3798.6
$ws.Range("I3").Value = 3359.8
$ws.Range("J3").Value = 4237.4
$ws.Range("K3").Value = 3359.8
$ws.Range("L3").Value = 4237.4
$ws.Range("M3").Value = -3245.8
$ws.Range("N3").Value = -4465.4
$ws.Range("H105").Value = 7237.4
$ws.Range("I105").Value = 13688.9
$ws.Range("J105").Value = 4011.65
$ws.Range("K105").Value = 13688.9
$ws.Range("L105").Value = 4011.65
$ws.Range("M105").Value = -11941.9
$ws.Range("N105").Value = -7505.65
$ws.Range("H107").Value = 4352.25
$ws.Range("I107").Value = 3227
$ws.Range("J107").Value = 5477.5
$ws.Range("K107").Value = 3227
$ws.Range("L107").Value = 5477.5
$ws.Range("M107").Value = -1307
$ws.Range("N107").Value = -9317.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2021.091
$ws.Range("I16").Value = 1802.375
$ws.Range("K16").Value = 1802.375
$ws.Range("M16").Value = -1515.375
$ws.Range("H31").Value = 7146631.5
$ws.Range("I31").Value = 2357.7917
$ws.Range("J31").Value = 22734138
$ws.Range("K31").Value = 2357.7917
$ws.Range("L31").Value = 22734138
$ws.Range("M31").Value = -2062.7917
$ws.Range("N31").Value = -22734728
$ws.Range("H34").Value = 7146631.5
$ws.Range("I34").Value = 2357.7917
$ws.Range("J34").Value = 22734138
$ws.Range("K34").Value = 2357.7917
$ws.Range("L34").Value = 22734138
$ws.Range("M34").Value = -2155.7917
$ws.Range("N34").Value = -22734542
$ws.Range("H58").Value = 2071.2666
$ws.Range("I58").Value = 1948.9
$ws.Range("K58").Value = 1948.9
$ws.Range("M58").Value = -1745.9
$ws.Range("H99").Value = 7856.4
$ws.Range("I99").Value = 4805.273
$ws.Range("J99").Value = 10253.714
$ws.Range("K99").Value = 4805.273
$ws.Range("L99").Value = 10253.714
$ws.Range("M99").Value = -3307.273
$ws.Range("N99").Value = -13249.714
$ws.Range("H107").Value = 1279
$ws.Range("I107").Value = 211
$ws.Range("J107").Value = 1991
$ws.Range("K107").Value = 211
$ws.Range("L107").Value = 1991
$ws.Range("M107").Value = 1709
$ws.Range("N107").Value = -5831
$ws.Range("H113").Value = 2021.091
$ws.Range("I113").Value = 1802.375
$ws.Range("K113").Value = 1802.375
$ws.Range("M113").Value = 367.625
$ws.Range("H126").Value = 7856.4
$ws.Range("I126").Value = 4805.273
$ws.Range("J126").Value = 10253.714
$ws.Range("K126").Value = 14415.819
$ws.Range("L126").Value = 30761.142
$ws.Range("M126").Value = -11945.819
$ws.Range("N126").Value = -35701.142
$ws.Range("H132").Value = 71611.2
$ws.Range("I132").Value = 93286.41
$ws.Range("K132").Value = 279859.23
$ws.Range("M132").Value = -277329.23
$ws.Range("H136").Value = 2071.2666
$ws.Range("I136").Value = 1948.9
$ws.Range("K136").Value = 5846.700000000001
$ws.Range("M136").Value = -3296.700000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 58906.06
$ws.Range("I2").Value = 67.53846
$ws.Range("J2").Value = 250131.25
$ws.Range("K2").Value = 405.23076
$ws.Range("L2").Value = 1500787.5
$ws.Range("M2").Value = -292.23076
$ws.Range("N2").Value = -1501013.5
$ws.Range("H15").Value = 297.66666
$ws.Range("I15").Value = 91
$ws.Range("J15").Value = 401
$ws.Range("K15").Value = 273
$ws.Range("L15").Value = 1203
$ws.Range("M15").Value = -133
$ws.Range("N15").Value = -1483
$ws.Range("H138").Value = 2359.5
$ws.Range("I138").Value = 1664.25
$ws.Range("K138").Value = 4992.75
$ws.Range("M138").Value = 147.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 35834.9
$ws.Range("I52").Value = 34900
$ws.Range("K52").Value = 34900
$ws.Range("M52").Value = -34641
$ws.Range("H107").Value = 704.7059
$ws.Range("J107").Value = 650
$ws.Range("L107").Value = 650
$ws.Range("N107").Value = -4490
$ws.Range("H122").Value = 38468068
$ws.Range("I122").Value = 7137.5
$ws.Range("J122").Value = 166671170
$ws.Range("K122").Value = 21412.5
$ws.Range("L122").Value = 500013510
$ws.Range("M122").Value = -18962.5
$ws.Range("N122").Value = -500018410

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H80").Value = 67654.8
$ws.Range("J80").Value = 60009.332
$ws.Range("L80").Value = 60009.332
$ws.Range("N80").Value = -62255.332
$ws.Range("H83").Value = 67654.8
$ws.Range("J83").Value = 60009.332
$ws.Range("L83").Value = 180027.996
$ws.Range("N83").Value = -191259.996
$ws.Range("H136").Value = 2005704.8
$ws.Range("I136").Value = 2503099.8
$ws.Range("K136").Value = 7509299.399999999
$ws.Range("M136").Value = -7506749.399999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 24664.334
$ws.Range("J45").Value = 24664.334
$ws.Range("L45").Value = 24664.334
$ws.Range("N45").Value = -25646.334
$ws.Range("H132").Value = 2440.8333
$ws.Range("I132").Value = 2365.5557
$ws.Range("K132").Value = 7096.6671
$ws.Range("M132").Value = -4566.6671
